# Generate Report for Handback
#
# Row 7 (430dd497-f5bf-49bb-bf6b-44afe203a876) on both the "zh-cn" and
# "de-de" sheets now has a completed handback: a "Latest Target File"
# hyperlink, a "Latest Handback File" name, a "Latest Handback DateTime"
# and an "Error Detail" message saying the handed-back file isn't built
# from the very latest source commit.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8219d336027fd890f5566d032e4a456526d49d0b/e2e/430dd497-f5bf-49bb-bf6b-44afe203a876.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/100941269c9d0945af2c124b32f7608c0ba40184/e2e/430dd497-f5bf-49bb-bf6b-44afe203a876.md."

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/100941269c9d0945af2c124b32f7608c0ba40184/e2e/430dd497-f5bf-49bb-bf6b-44afe203a876.md",
    "",
    "",
    "430dd497-f5bf-49bb-bf6b-44afe203a876.md"
)
$wsZh.Range("I7").Font.Underline = $true
$wsZh.Range("I7").Font.Color = 15570276

$wsZh.Range("J7").Value = "430dd497-f5bf-49bb-bf6b-44afe203a876.ff6d6c96366b78ca71cd2e31586ac7d7c852de6f.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-03 06:59:06"
$wsZh.Range("P7").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/100941269c9d0945af2c124b32f7608c0ba40184/e2e/430dd497-f5bf-49bb-bf6b-44afe203a876.md",
    "",
    "",
    "430dd497-f5bf-49bb-bf6b-44afe203a876.md"
)
$wsDe.Range("I7").Font.Underline = $true
$wsDe.Range("I7").Font.Color = 15570276

$wsDe.Range("J7").Value = "430dd497-f5bf-49bb-bf6b-44afe203a876.ff6d6c96366b78ca71cd2e31586ac7d7c852de6f.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-03 06:59:14"
$wsDe.Range("P7").Value = $errorDetail
